{"js": "// Convert a handful of table header labels in the \"vacant sale\" template\n// from sentence case to Title Case, and append \" (Days)\" to the\n// \"Time on market\" label so it reads \"Time on Market (Days)\".\n//\n// The run formatting (font, size, bCs, etc.) on every label is identical\n// before and after the edit, so a straightforward text replace on each\n// matched range preserves the existing character formatting exactly.\n\nconst replacements = [\n  [\"Legal description\", \"Legal Description\"],\n  [\"Sale date\", \"Sale Date\"],\n  [\"Property rights\", \"Property Rights\"],\n  [\"Conditions of sale\", \"Conditions of Sale\"],\n  [\"Time on market\", \"Time on Market (Days)\"],\n  [\"Sale price/SF\", \"Sale Price/SF\"],\n  [\"Land size (square feet)\", \"Land Size (Square Feet)\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: '\" + oldText + \"'\");\n  }\n\n  for (const found of results.items) {\n    found.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Convert a handful of table header labels in the \"vacant sale\" template\n# from sentence case to Title Case, and append \" (Days)\" to the\n# \"Time on market\" label so it reads \"Time on Market (Days)\".\n#\n# The run formatting (font, size, bCs, etc.) on every label is identical\n# before and after the edit, so a plain Find/Replace on each label's text\n# preserves the existing character formatting exactly.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"Legal description\", \"Legal Description\"),\n  @(\"Sale date\", \"Sale Date\"),\n  @(\"Property rights\", \"Property Rights\"),\n  @(\"Conditions of sale\", \"Conditions of Sale\"),\n  @(\"Time on market\", \"Time on Market (Days)\"),\n  @(\"Sale price/SF\", \"Sale Price/SF\"),\n  @(\"Land size (square feet)\", \"Land Size (Square Feet)\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $rng.Find.Found) {\n    throw \"Could not find text to replace: '$oldText'\"\n  }\n}\n"}
